# "Common: Welcome gifts are taken from pricelist"
#
# The "prices" sheet gains two new pricelist rows for the welcome-gift
# tariffs (user/root), reusing the existing "JX-BMD9-GYJXO9" tariff code,
# and the whole data block (A2:C4) picks up the "import" look (10pt font
# + wrapped text) used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prices")

# Row 2 stays the same data (lab.liquid.create @ 10), it just picks up the
# new formatting applied to the whole A2:C4 block below.

# Row 3: welcome-gift.user @ 500
$ws.Range("A3").Value = "JX-BMD9-GYJXO9"
$ws.Range("B3").Value = "welcome-gift.user"
$ws.Range("C3").Value = 500

# Row 4: welcome-gift.root @ 10000 (was a different tariff/price before)
$ws.Range("A4").Value = "JX-BMD9-GYJXO9"
$ws.Range("B4").Value = "welcome-gift.root"
$ws.Range("C4").Value = 10000

# Apply the "import" look (10pt font, wrapped text) to the data rows.
$data = $ws.Range("A2:C4")
$data.Font.Size = 10
$data.WrapText = $true

# Move the active selection to C4 (last edited cell).
$ws.Range("C4").Select()
